$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 14: "DRG 1975 IV" was a typo, should be "DRG 1975 VI"
#     (matches the Kürzel "DRG1975VI" already in column A)
$ws.Range("B14").Value = "DRG 1975 VI"

# --- Re-measured header-row heights (rows 1-9 carry custom, non-default
#     formatting; everything below uses the sheet default of 15.75 and is
#     untouched)
$ws.Rows("1").RowHeight = 13
$ws.Rows("2").RowHeight = 13
$ws.Rows("3").RowHeight = 25.5
$ws.Rows("4").RowHeight = 25
$ws.Rows("5").RowHeight = 25
$ws.Rows("6").RowHeight = 25
$ws.Rows("7").RowHeight = 12.5
$ws.Rows("8").RowHeight = 12.5
$ws.Rows("9").RowHeight = 25

# --- Append 10 new literature rows (100-109), copying the formatting of the
#     last existing data row (99) so the new rows pick up the same cell
#     styles, then inserting them (shifted down) right after row 99.
$ws.Rows("99").Copy()
$ws.Rows("100:109").Insert(-4121)

# Row 100: Schabus 1971 I
$ws.Range("A100").Value = "Schabus1971I"
$ws.Range("B100").Value = "Schabus 1971 I"
$ws.Range("C100").Value = 1971
$ws.Range("D100").Value = "Dialektgeographie des Lesachtals (Kärnten) - 1"
$ws.Range("E100").Value = "Schabus"
$ws.Range("F100").Value = "Wilfried"

# Row 101: Schabus 1971 II
$ws.Range("A101").Value = "Schabus1971II"
$ws.Range("B101").Value = "Schabus 1971 II"
$ws.Range("C101").Value = 1971
$ws.Range("D101").Value = "Dialektgeographie des Lesachtals (Kärnten) - 2"
$ws.Range("E101").Value = "Schabus"
$ws.Range("F101").Value = "Wilfried"

# Row 102: Pirchegger 1927
$ws.Range("A102").Value = "Pirchegger1927"
$ws.Range("B102").Value = "Pirchegger 1927"
$ws.Range("C102").Value = 1927
$ws.Range("D102").Value = "Die slavischen Ortsnamen im Mürzgebiet"
$ws.Range("E102").Value = "Pirchegger"
$ws.Range("F102").Value = "Simon"

# Row 103: Kronen Zeitung
$ws.Range("A103").Value = "Kronen_Ztg"
$ws.Range("B103").Value = "Kronen Zeitung"
$ws.Range("C103").Value = 1980
$ws.Range("D103").Value = "Kronen Zeitung 28.9.1980"

# Row 104: Pailler 1883 II
$ws.Range("A104").Value = "Pailler1883II"
$ws.Range("B104").Value = "Pailler 1883 II"
$ws.Range("C104").Value = 1883
$ws.Range("D104").Value = "Krippenspiele aus Oberösterreich und Tirol"
$ws.Range("E104").Value = "Pailler"
$ws.Range("F104").Value = "Wilhelm"

# Row 105: Pailler 1881 I
$ws.Range("A105").Value = "Pailler1881I"
$ws.Range("B105").Value = "Pailler 1881 I"
$ws.Range("C105").Value = 1881
$ws.Range("D105").Value = "Weihnachtlieder aus Oberösterreich"
$ws.Range("E105").Value = "Pailler"
$ws.Range("F105").Value = "Wilhelm"

# Row 106: Schatzdorfer 1949
$ws.Range("A106").Value = "Schatzdorfer1949"
$ws.Range("B106").Value = "Schatzdorfer 1949"
$ws.Range("C106").Value = 1949
$ws.Range("D106").Value = "Spatzngsang und Spinnáwittn"
$ws.Range("E106").Value = "Schatzdorfer"
$ws.Range("F106").Value = "Hans"

# Row 107: Goldbacher 1904
$ws.Range("A107").Value = "Goldbacher1904"
$ws.Range("B107").Value = "Goldbacher 1904"
$ws.Range("C107").Value = 1904
$ws.Range("D107").Value = "Gmüatlichö Sach'n"
$ws.Range("E107").Value = "Goldbacher"
$ws.Range("F107").Value = "Gregor"

# Row 108: Der Schlern 1920
$ws.Range("A108").Value = "DerSchlern1920"
$ws.Range("B108").Value = "Der Schlern 1920"
$ws.Range("C108").Value = 1920
$ws.Range("D108").Value = "Bozner Halbmonatsschrift"

# Row 109: Rudl 1920
$ws.Range("A109").Value = "Rudl1920"
$ws.Range("B109").Value = "Rudl 1920"
$ws.Range("C109").Value = 1920
$ws.Range("D109").Value = "Der Hiesl ban Zonndoktr. Der Schlern. Ausgabe vom 1.6.1920"
$ws.Range("E109").Value = "Rudl"
$ws.Range("F109").Value = "Otto  "

# --- Grow the Tabelle1 table/ListObject to cover the newly appended rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:H109"))

# --- Update the view: scroll so the new rows are visible and select C112
#     (mirrors the author's on-screen state when the file was saved)
$win = $excel.ActiveWindow
$win.ScrollRow = 95
$win.ScrollColumn = 1
$null = $ws.Range("C112").Select()
